$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filter2-BOM")

# Fill in column C for the BOM rows so the "on hand" quantity matches the
# required quantity, which drives the shared formula in column E to 0.
$ws.Range("C10").Value = 4
$ws.Range("C24").Value = 2
$ws.Range("C25").Value = 2
$ws.Range("C26").Value = 2

# Move the active selection to C11 (was D28)
$ws.Range("C11").Select()
